# Scheduled market-data refresh: update price/profit columns (H-N) across all Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 734.9091
$ws.Range("I15").Value = 734.9091
$ws.Range("K15").Value = 2204.7273
$ws.Range("M15").Value = -2035.7273
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 2499.5
$ws.Range("J29").Value = 2499.5
$ws.Range("L29").Value = 7498.5
$ws.Range("N29").Value = -8060.5
$ws.Range("H98").Value = 930.5
$ws.Range("J98").Value = 995
$ws.Range("L98").Value = 995
$ws.Range("N98").Value = -3991
$ws.Range("H122").Value = 930.5
$ws.Range("J122").Value = 995
$ws.Range("L122").Value = 2985
$ws.Range("N122").Value = -7885
$ws.Range("H132").Value = 13518.333
$ws.Range("I132").Value = 13518.333
$ws.Range("K132").Value = 40554.999
$ws.Range("M132").Value = -38024.999
$ws.Range("H137").Value = 2979.6667
$ws.Range("I137").Value = 2283.3333
$ws.Range("J137").Value = 3676
$ws.Range("K137").Value = 6849.999899999999
$ws.Range("L137").Value = 11028
$ws.Range("M137").Value = -4299.999899999999
$ws.Range("N137").Value = -16128
$ws.Range("H138").Value = 8294.1
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 8294.1
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 24882.3
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -35162.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 347.75
$ws.Range("I26").Value = 347.75
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 347.75
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -17.75
$ws.Range("N26").ClearContents()
$ws.Range("H41").Value = 13820
$ws.Range("I41").Value = 13820
$ws.Range("K41").Value = 13820
$ws.Range("M41").Value = -13406
$ws.Range("H102").Value = 1659.2858
$ws.Range("I102").Value = 1659.2858
$ws.Range("K102").Value = 1659.2858
$ws.Range("M102").Value = -37.28580000000011
$ws.Range("H132").Value = 6005.5454
$ws.Range("I132").Value = 5848.2
$ws.Range("K132").Value = 17544.6
$ws.Range("M132").Value = -15014.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 6200
$ws.Range("J49").Value = 6200
$ws.Range("L49").Value = 6200
$ws.Range("N49").Value = -6678
$ws.Range("H82").Value = 20841.25
$ws.Range("J82").Value = 130000
$ws.Range("L82").Value = 130000
$ws.Range("N82").Value = -130766
$ws.Range("H85").Value = 20841.25
$ws.Range("J85").Value = 130000
$ws.Range("L85").Value = 130000
$ws.Range("N85").Value = -132652
$ws.Range("H105").Value = 1103.3334
$ws.Range("I105").Value = 1099.5
$ws.Range("K105").Value = 1099.5
$ws.Range("M105").Value = 647.5
$ws.Range("H134").Value = 2800
$ws.Range("J134").Value = 4900
$ws.Range("L134").Value = 14700
$ws.Range("N134").Value = -19770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 286.6
$ws.Range("I32").Value = 286.6
$ws.Range("K32").Value = 286.6
$ws.Range("M32").Value = 29.39999999999998
$ws.Range("H35").Value = 2669.625
$ws.Range("I35").Value = 896.75
$ws.Range("J35").Value = 4442.5
$ws.Range("K35").Value = 896.75
$ws.Range("L35").Value = 4442.5
$ws.Range("M35").Value = -602.75
$ws.Range("N35").Value = -5030.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 908
$ws.Range("I80").Value = 897.5
$ws.Range("J80").Value = 950
$ws.Range("K80").Value = 2692.5
$ws.Range("L80").Value = 2850
$ws.Range("M80").Value = -1756.5
$ws.Range("N80").Value = -4722
$ws.Range("H83").Value = 908
$ws.Range("I83").Value = 897.5
$ws.Range("J83").Value = 950
$ws.Range("K83").Value = 8077.5
$ws.Range("L83").Value = 8550
$ws.Range("M83").Value = -3397.5
$ws.Range("N83").Value = -17910
$ws.Range("H92").Value = 197.33333
$ws.Range("I92").Value = 392
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 1176
$ws.Range("L92").Value = 300
$ws.Range("M92").Value = 72
$ws.Range("N92").Value = -2796
$ws.Range("H116").Value = 2581.7273
$ws.Range("I116").Value = 1700
$ws.Range("J116").Value = 2777.6667
$ws.Range("K116").Value = 5100
$ws.Range("L116").Value = 8333.000100000001
$ws.Range("M116").Value = -1658
$ws.Range("N116").Value = -15217.0001
$ws.Range("H117").Value = 5279.8
$ws.Range("I117").Value = 705.6667
$ws.Range("J117").Value = 12141
$ws.Range("K117").Value = 2117.0001
$ws.Range("L117").Value = 36423
$ws.Range("M117").Value = 1324.9999
$ws.Range("N117").Value = -43307
$ws.Range("H129").Value = 12285
$ws.Range("J129").Value = 19921.666
$ws.Range("L129").Value = 59764.99800000001
$ws.Range("N129").Value = -69764.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17331.666
$ws.Range("I70").Value = 3500
$ws.Range("K70").Value = 3500
$ws.Range("M70").Value = -3230
$ws.Range("H73").Value = 17331.666
$ws.Range("I73").Value = 3500
$ws.Range("K73").Value = 3500
$ws.Range("M73").Value = -2564
$ws.Range("H102").Value = 768.4286
$ws.Range("I102").Value = 675.8
$ws.Range("K102").Value = 675.8
$ws.Range("M102").Value = 946.2
$ws.Range("H122").Value = 1151.7142
$ws.Range("I122").Value = 1151.7142
$ws.Range("K122").Value = 3455.1426
$ws.Range("M122").Value = -1005.1426
$ws.Range("H132").Value = 1890.3334
$ws.Range("J132").Value = 1883
$ws.Range("L132").Value = 5649
$ws.Range("N132").Value = -10709

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 1500
$ws.Range("I32").Value = 1500
$ws.Range("K32").Value = 1500
$ws.Range("M32").Value = -1183
$ws.Range("H132").Value = 4777.4
$ws.Range("I132").Value = 4777.4
$ws.Range("K132").Value = 14332.2
$ws.Range("M132").Value = -11802.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 42999.75
$ws.Range("I58").Value = 42999.75
$ws.Range("K58").Value = 42999.75
$ws.Range("M58").Value = -42691.75
$ws.Range("H62").Value = 2666.6667
$ws.Range("I62").Value = 2666.6667
$ws.Range("K62").Value = 2666.6667
$ws.Range("M62").Value = -2042.6667
$ws.Range("H65").Value = 2666.6667
$ws.Range("I65").Value = 2666.6667
$ws.Range("K65").Value = 13333.3335
$ws.Range("M65").Value = -10213.3335
$ws.Range("H132").Value = 1782.4546
$ws.Range("I132").Value = 1560.8
$ws.Range("K132").Value = 4682.4
$ws.Range("M132").Value = -2152.4
$ws.Range("H136").Value = 2489.3076
$ws.Range("I136").Value = 2566.0908
$ws.Range("J136").Value = 2067
$ws.Range("K136").Value = 7698.2724
$ws.Range("L136").Value = 6201
$ws.Range("M136").Value = -5148.2724
$ws.Range("N136").Value = -11301
